$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version 0.1.6 -> 0.1.7
$ws1.Range("B3").Value = "0.1.7"

# Status active -> draft
$ws1.Range("B6").Value = "draft"

# Date updated
$ws1.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Contact row 10 detail text updated (publisher contact)
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Contact row 11 detail text updated (individual contact)
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# A new row 12 ("Jurisdiction") is inserted, pushing Description/Purpose/Copyright/
# Immutable down by one row. Avoid Rows.Insert() (it mints a spurious extra cell
# style in this runtime) - shift the values down manually instead, bottom-up so
# nothing is clobbered before it is read, and copy the formatting of the existing
# data rows onto the newly-used row 16 so every row keeps the shared style "s=2".

$ws1.Range("A16").Value = "Immutable"
$ws1.Range("B16").Value = "BooleanType[null]"
$ws1.Range("A15:B15").Copy()
$ws1.Range("A16:B16").PasteSpecial(-4122) # xlPasteFormats

$ws1.Range("A15").Value = "Copyright"
$ws1.Range("B15").Value = ""

$ws1.Range("A14").Value = "Purpose"
$ws1.Range("B14").Value = ""

$ws1.Range("A13").Value = "Description"
$ws1.Range("B13").Value = "SNOMED: Disorders of hematopoietic structure"

$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""

Write-Output "done"
